# Update scraped batting performance counts (ODI/D column and T20/E column)
# for various players in the "all_formats_raw" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_formats_raw")

$updates = @{
    "D2"  = 26
    "D7"  = 140
    "E7"  = 80
    "E13" = 15
    "E17" = 51
    "D18" = 9
    "D19" = 36
    "E19" = 41
    "D22" = 7
    "D25" = 50
    "E25" = 34
    "D26" = 155
    "E26" = 114
    "D29" = 45
    "E29" = 37
    "D30" = 21
    "E30" = 31
    "E32" = 56
    "D36" = 89
    "E36" = 56
    "E38" = 29
    "D39" = 44
    "D45" = 45
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
